$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet3")

# --- Rebuild the little QTL R-square block as a proper 3-column table ---
# Old layout (row 1: headers Del/Nig/Dos/Combined in B1:E1, row 2: R-square
# label in A2 with values in B2:E2) becomes a transposed table with headers
# Models / # of QTLs / Variance explained (%) in A1:C1 and one row per model.

# Apply the new MyriadPro-Regular font across the old A1:E2 footprint first
$ws.Range("A1:E2").Font.Name = "MyriadPro-Regular"

# Header row (bold). Variance explained (%) is written before # of QTLs so
# the shared-string table picks up the same ordering as the target file.
$ws.Range("A1").Value = "Models"
$ws.Range("C1").Value = "Variance explained (%)"
$ws.Range("B1").Value = "# of QTLs"
$ws.Range("A1:C1").Font.Bold = $true

# Data rows
$ws.Range("A2").Value = "Del"
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 2.34

$ws.Range("A3").Value = "Nig"
$ws.Range("B3").Value = 3
$ws.Range("C3").Value = 13.88

$ws.Range("A4").Value = "Dos"
$ws.Range("B4").Value = 10
$ws.Range("C4").Value = 15.19

$ws.Range("A5").Value = "Combined"
$ws.Range("B5").Value = 7
$ws.Range("C5").Value = 88.22

# Clear the leftover old D/E values beyond the new 3-column table (their
# formatting was already set above when A1:E2 got the font change)
$ws.Range("D1:E2").ClearContents()

# Extend the matching font down the rest of column A and C
$ws.Range("A3").Font.Name = "MyriadPro-Regular"
$ws.Range("A5").Font.Name = "MyriadPro-Regular"
$ws.Range("C3:C5").Font.Name = "MyriadPro-Regular"

# Column widths (close to the author's best-fit sizing for the new content)
$ws.Columns.Item(1).ColumnWidth = 11
$ws.Columns.Item(2).ColumnWidth = 11.8
$ws.Columns.Item(3).ColumnWidth = 20

# Selection to match the saved view state
$ws.Range("A1:C5").Select()
